$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# --- Rebuild the "Código" column as a new table column named "Código " ---
# (the author deleted the old "Código" column and re-added a fresh one,
#  which is why the table XML gets a brand-new column id/uid downstream)
$oldCol = $tbl.ListColumns.Item("Código")
$oldCol.Delete()
$newCol = $tbl.ListColumns.Add()

# Header text (note trailing space, matches the source workbook)
$ws.Range("F1").Value = "Código "

# Column width for the new "Código " column (engine stores width + 5/6,
# so back the requested width off by that padding to land on 18 exactly)
$ws.Columns.Item(6).ColumnWidth = 17.1666667

# --- Fill in the product codes for every data row ---
$codes = @{
    2  = "FX1";  3  = "FX2";  4  = "FX3";  5  = "FX4";  6  = "FX5";
    7  = "FX6";  8  = "FX7";  9  = "FX8";  10 = "FX9";  11 = "FX10";
    12 = "FX11"; 13 = "FX12"; 14 = "FX13";
    15 = "FS1";  16 = "FS2";  17 = "FS3";  18 = "FS4";  19 = "FS5";  20 = "FS6";
    21 = "FS7";  22 = "FS8";  23 = "FS9";  24 = "FS10"; 25 = "FS11"; 26 = "FS12";
    27 = "FS13"; 28 = "FS14";
    29 = "FI1";  30 = "FI2";  31 = "FI3";  32 = "FI4"
}

foreach ($row in $codes.Keys) {
    $ws.Cells.Item($row, 6).Value = $codes[$row]
}

# --- Update the view: scroll/selection as left by the author ---
$ws.Range("G33").Select() | Out-Null
